# Update 农用柴油和农药使用量.xlsx : drop years 2000-2009 and keep only
# 2010-2019 (shifted up to occupy rows 2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..11 (previously rows 12..21: years 2010-2019)
$data = @(
    @("2010年", 2023.1181,     175.8219272),
    @("2011年", 2057.4404015,  178.7001839),
    @("2012年", 2107.6484,     180.6057204),
    @("2013年", 2154.943653,   180.770375),
    @("2014年", 2176.317746,   180.33147194),
    @("2015年", 2197.6810583,  178.29693251),
    @("2016年", 2117.0837779,  174.04585952),
    @("2017年", 2095.1126237,  165.506603081),
    @("2018年", 2003.391689,   150.355276843),
    @("2019年", 1934,          139.1747)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $ws.Cells.Item($rowIndex, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowIndex, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowIndex, 3).Value = $data[$i][2]
}

# Remove the now-duplicated trailing rows (old rows 12-21), leaving the
# sheet with just rows 1-11.
$ws.Rows("12:21").Delete()
